$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates (report volume/issue number and the week's date
# range). These are rich-text shared strings; use Characters() so we
# only touch the specific substring that changed, mirroring how this
# would be hand-edited in Excel.
# ---------------------------------------------------------------------

# A8: "Volume 30   Number  37" -> "Volume 30   Number  38"
$volRange = $ws.Range("A8")
$volText = $volRange.Characters(1, 200).Text
$volIdx = $volText.LastIndexOf("37") + 1
$volRange.Characters($volIdx, 2).Text = "38"

# C9: "Report Covering the Week  9/11/2023  Through  9/17/2023"
#  -> "Report Covering the Week  9/18/2023  Through  9/24/2023"
$weekRange = $ws.Range("C9")
$weekText = $weekRange.Characters(1, 200).Text
$startIdx = $weekText.IndexOf("9/11/2023") + 1
$weekRange.Characters($startIdx, 9).Text = "9/18/2023"
$weekText2 = $weekRange.Characters(1, 200).Text
$endIdx = $weekText2.IndexOf("9/17/2023") + 1
$weekRange.Characters($endIdx, 9).Text = "9/24/2023"

# ---------------------------------------------------------------------
# Crime-stat table updates (rows 16-29). New weekly crime figures were
# collected, which ripple through the weekly/28-day/YTD counts and all
# of their derived percent-change columns.
# ---------------------------------------------------------------------

# A few cells flip from the "not applicable" placeholder (shared text
# "0" / "***.*") to real numbers now that data exists for them; give
# them the same number format as their numeric neighbors before writing
# the value so they pick up the right style (counts vs. percentages).
$countFormat = $ws.Range("D16").NumberFormat
$pctFormat = $ws.Range("E16").NumberFormat

$ws.Range("D22").NumberFormat = $countFormat
$ws.Range("E22").NumberFormat = $pctFormat
$ws.Range("C26").NumberFormat = $countFormat
$ws.Range("C27").NumberFormat = $countFormat
$ws.Range("D27").NumberFormat = $countFormat
$ws.Range("E27").NumberFormat = $pctFormat

# Row 16 - Robbery
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 156
$ws.Range("J16").Value = 168
$ws.Range("K16").Value = -7.142857142857
$ws.Range("L16").Value = 13.043478260869
$ws.Range("M16").Value = -16.577540106951
$ws.Range("N16").Value = -75.04

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 40
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 17.391304347826
$ws.Range("I17").Value = 263
$ws.Range("J17").Value = 257
$ws.Range("K17").Value = 2.334630350194
$ws.Range("L17").Value = 10.041841004184
$ws.Range("M17").Value = 152.884615384615
$ws.Range("N17").Value = -13.201320132013

# Row 18 - Burglary
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -37.5
$ws.Range("I18").Value = 103
$ws.Range("J18").Value = 104
$ws.Range("K18").Value = -0.961538461538
$ws.Range("L18").Value = 49.275362318840
$ws.Range("M18").Value = -52.534562211981
$ws.Range("N18").Value = -87.967289719626

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 15.384615384615
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = 33.333333333333
$ws.Range("I19").Value = 464
$ws.Range("J19").Value = 488
$ws.Range("K19").Value = -4.918032786885
$ws.Range("L19").Value = 32.951289398280
$ws.Range("M19").Value = 90.946502057613
$ws.Range("N19").Value = 9.692671394799

# Row 20 - G.L.A.
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 181
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = 13.125
$ws.Range("L20").Value = 41.40625
$ws.Range("M20").Value = -17.351598173516
$ws.Range("N20").Value = -92.366090257275

# Row 21 - TOTAL
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 17.241379310344
$ws.Range("F21").Value = 128
$ws.Range("G21").Value = 109
$ws.Range("H21").Value = 17.431192660550
$ws.Range("I21").Value = 1183
$ws.Range("J21").Value = 1199
$ws.Range("K21").Value = -1.334445371142
$ws.Range("L21").Value = 26.254002134471
$ws.Range("M21").Value = 19.615773508594
$ws.Range("N21").Value = -74.393939393939

# Row 22 - Transit
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 19
$ws.Range("K22").Value = -10.526315789473

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 40.909090909090
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 116
$ws.Range("H24").Value = -17.241379310344
$ws.Range("I24").Value = 942
$ws.Range("J24").Value = 1041
$ws.Range("K24").Value = -9.510086455331
$ws.Range("L24").Value = 45.820433436532
$ws.Range("M24").Value = 99.576271186440

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -9.090909090909
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = -5.714285714285
$ws.Range("I25").Value = 397
$ws.Range("J25").Value = 392
$ws.Range("K25").Value = 1.275510204081
$ws.Range("L25").Value = 15.743440233236
$ws.Range("M25").Value = 3.116883116883

# Row 26 - UCR Rape*
$ws.Range("C26").Value = 1
$ws.Range("I26").Value = 25
$ws.Range("K26").Value = -10.714285714285
$ws.Range("L26").Value = 66.666666666666

# Row 27 - Other Sex Crimes
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -20
$ws.Range("I27").Value = 38
$ws.Range("J27").Value = 42
$ws.Range("K27").Value = -9.523809523809
$ws.Range("L27").Value = 22.580645161290

# Row 28 - Shooting Vic.
$ws.Range("L28").Value = -68.75

# Row 29 - Shooting Inc.
$ws.Range("L29").Value = -75
